$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per row based on the re-aligned weekly data (Fruta / hortaliza, semanal).
$ws.Range("D2").Value = 44214
$ws.Range("M2").Value = 48
$ws.Range("R2").Value = "Provincia de Linares"
$ws.Range("D3").Value = 44592
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 4000
$ws.Range("D4").Value = 44614
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000
$ws.Range("D5").Value = 44628
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("S5").Value = 3000
$ws.Range("D6").Value = 44582
$ws.Range("M6").Value = 150
$ws.Range("O6").Value = 6500
$ws.Range("P6").Value = 6233
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 3116
$ws.Range("D7").Value = 44588
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 6500
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 6750
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 3375
$ws.Range("D8").Value = 44585
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 6500
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 6750
$ws.Range("S8").Value = 3375
$ws.Range("D9").Value = 44627
$ws.Range("M9").Value = 45
$ws.Range("D10").Value = 44586
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 7000
$ws.Range("P10").Value = 7000
$ws.Range("S10").Value = 3500
$ws.Range("D11").Value = 44587
$ws.Range("M11").Value = 165
$ws.Range("N11").Value = 6500
$ws.Range("O11").Value = 7000
$ws.Range("P11").Value = 6742
$ws.Range("R11").Value = "Provincia de Linares"
$ws.Range("S11").Value = 3371
$ws.Range("D12").Value = 44209
$ws.Range("M12").Value = 58
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("D13").Value = 44211
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 6000
$ws.Range("O13").Value = 6000
$ws.Range("P13").Value = 6000
$ws.Range("S13").Value = 3000
$ws.Range("D14").Value = 44589
$ws.Range("M14").Value = 60
$ws.Range("R14").Value = "Provincia de Curicó"
$ws.Range("D15").Value = 44606
$ws.Range("M15").Value = 45
$ws.Range("R15").Value = "Provincia de Linares"
